$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (columns B:E)
$ws.Range("B2").Value = 8.7353400442822213
$ws.Range("C2").Value = 5.7490026772609859
$ws.Range("D2").Value = 9.4831743166300608
$ws.Range("E2").Value = 7.8741025572924395

# Row 3 data values (columns B:E)
$ws.Range("B3").Value = 6.0993965164398682
$ws.Range("C3").Value = 7.1870464587086405
$ws.Range("D3").Value = 5.6930233603028739
$ws.Range("E3").Value = 8.3593872193739411

# Update the selection to reflect the new active range (B1:E3) as in the saved file
$ws.Range("B1:E3").Select()
